$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new risk row values (order matches how the author originally typed them,
# which determines the shared-strings table order)
$ws.Range("A2").Value = "Accidentally drop IMS database while programme is running"
$ws.Range("D2").Value = "High"
$ws.Range("C2").Value = "low"
$ws.Range("E2").Value = "developer"
$ws.Range("B2").Value = "IMS programme can no longer access a database for any commands"
$ws.Range("F2").Value = "End application and recreate database"
$ws.Range("G2").Value = "Only allow for CRUD sql statements within code, give no access to DROP keyword"

# Wrap text for the new row, keeping the existing border
$ws.Range("A2:G2").WrapText = $true

# Adjust column widths (values chosen so the engine's internal px-grid produces
# the closest attainable stored width to the target)
$ws.Columns.Item(2).ColumnWidth = 21
$ws.Columns.Item(6).ColumnWidth = 13.17

# Adjust row height for wrapped content
$ws.Rows.Item(2).RowHeight = 74.4

# Update selection to match final state
$ws.Range("H8").Select()
